# Apply the "element_isenabled,element_isdisabled are included in the list" edit.
$wb = $excel.ActiveWorkbook

$objectsWs = $wb.Worksheets.Item("Objects")
$toolbarWs = $wb.Worksheets.Item("Objects_Toolbar")

# --- 1. Objects sheet: remove "getCurrentFrameName" (row 31) and append new rows ---
# Current (before) layout:
#   A29 = switchToFrameFromDefault
#   A30 = switchToSingleFrame
#   A31 = getCurrentFrameName
# Target layout:
#   A29 = switchToFrameFromDefault   (unchanged)
#   A30 = switchToSingleFrame        (unchanged)
#   A31 = textbox_autosuggest_choice (replaces getCurrentFrameName)
#   A32 = textbox_autosuggest_browse
#   A33 = selectitembytextfromlist
#   A34 = element_enable
#   A35 = element_disable
#   A36 = element_displayed
#   A37 = element_notdisplayed
$objectsWs.Range("A31").Value = "textbox_autosuggest_choice"
$objectsWs.Range("A32").Value = "textbox_autosuggest_browse"
$objectsWs.Range("A33").Value = "selectitembytextfromlist"
$objectsWs.Range("A34").Value = "element_enable"
$objectsWs.Range("A35").Value = "element_disable"
$objectsWs.Range("A36").Value = "element_displayed"
$objectsWs.Range("A37").Value = "element_notdisplayed"

# --- 2. Update sheet view state on Objects sheet (scroll position / selection) ---
$objectsWs.Activate()
$excel.ActiveWindow.ScrollRow = 19
$objectsWs.Range("A39").Select()

# --- 3. Update sheet view/selection on Objects_Toolbar sheet (drop explicit selection) ---
$toolbarWs.Activate()
$toolbarWs.Range("A1").Select()

# --- 4. Update the x14 (Excel 2010+) list data validation on Objects_Toolbar!D column ---
# Before: one rule, sqref D3:D4, source Objects!$A$2:$A$28
# After:  two rules
#   rule A: sqref D4,        source Objects!$A$2:$A$50
#   rule B: sqref D3 D5:D22, source Objects!$A$2:$A$50
$toolbarWs.Range("D3:D4").Validation.Delete()

$toolbarWs.Range("D4").Validation.Add(3, 1, 1, "=Objects!`$A`$2:`$A`$50")
$toolbarWs.Range("D4").Validation.InputMessage = ""
$toolbarWs.Range("D4").Validation.ErrorMessage = ""
$toolbarWs.Range("D4").Validation.ShowInput = $true
$toolbarWs.Range("D4").Validation.ShowError = $true

$toolbarWs.Range("D3,D5:D22").Validation.Add(3, 1, 1, "=Objects!`$A`$2:`$A`$50")
$toolbarWs.Range("D3,D5:D22").Validation.InputMessage = ""
$toolbarWs.Range("D3,D5:D22").Validation.ErrorMessage = ""
$toolbarWs.Range("D3,D5:D22").Validation.ShowInput = $true
$toolbarWs.Range("D3,D5:D22").Validation.ShowError = $true

# --- 5. Workbook view: hide "Objects" sheet, keep "Objects_Toolbar" active, and
#        set firstSheet so the tab scroll starts at the second sheet ---
$objectsWs.Visible = $false
$toolbarWs.Activate()
$excel.ActiveWindow.DisplayWorkbookTabs = $true
$wb.Windows.Item(1).ScrollWorkbookTabs(1)
